$wb = $excel.ActiveWorkbook

# --- Flujo sheet (sheet5) edits ---
$ws = $wb.Worksheets.Item("Flujo")

# Clear row 4 entirely (old workflow-3 "3,1,2,1" row removed)
$ws.Rows.Item(4).Clear()

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1
$ws.Range("F2").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G2").Formula = '=A2&","&B2&","&C2&","&D2&")"'
$ws.Range("I2").Formula = '=F2&G2&";"'
$ws.Range("J2").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (5,1,2,1);'

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G3").Formula = '=A3&","&B3&","&C3&","&D3&")"'
$ws.Range("I3").Formula = '=F3&G3&";"'
$ws.Range("J3").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (5,2,3,2);'

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("F5").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G5").Formula = '=A5&","&B5&","&C5&","&D5&")"'
$ws.Range("I5").Formula = '=F5&G5&";"'
$ws.Range("J5").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,1,2,1);'

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 2
$ws.Range("F6").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G6").Formula = '=A6&","&B6&","&C6&","&D6&")"'
$ws.Range("I6").Formula = '=F6&G6&";"'
$ws.Range("J6").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,2,3,2);'

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 3
$ws.Range("F7").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G7").Formula = '=A7&","&B7&","&C7&","&D7&")"'
$ws.Range("I7").Formula = '=F7&G7&";"'
$ws.Range("J7").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,3,4,3);'

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 4
$ws.Range("F8").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G8").Formula = '=A8&","&B8&","&C8&","&D8&")"'
$ws.Range("I8").Formula = '=F8&G8&";"'
$ws.Range("J8").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,4,5,4);'

$ws.Range("A9").Value = 2
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 5
$ws.Range("F9").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G9").Formula = '=A9&","&B9&","&C9&","&D9&")"'
$ws.Range("I9").Formula = '=F9&G9&";"'
$ws.Range("J9").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,4,6,5);'

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 6
$ws.Range("F10").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G10").Formula = '=A10&","&B10&","&C10&","&D10&")"'
$ws.Range("I10").Formula = '=F10&G10&";"'
$ws.Range("J10").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,2,4,6);'

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 7
$ws.Range("F11").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G11").Formula = '=A11&","&B11&","&C11&","&D11&")"'
$ws.Range("I11").Formula = '=F11&G11&";"'
$ws.Range("J11").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,1,7,7);'

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 8
$ws.Range("F12").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G12").Formula = '=A12&","&B12&","&C12&","&D12&")"'
$ws.Range("I12").Formula = '=F12&G12&";"'
$ws.Range("J12").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,5,7,8);'

$ws.Range("A13").Value = 2
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 9
$ws.Range("F13").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G13").Formula = '=A13&","&B13&","&C13&","&D13&")"'
$ws.Range("I13").Formula = '=F13&G13&";"'
$ws.Range("J13").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,5,4,9);'

$ws.Range("A14").Value = 2
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 10
$ws.Range("F14").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G14").Formula = '=A14&","&B14&","&C14&","&D14&")"'
$ws.Range("I14").Formula = '=F14&G14&";"'
$ws.Range("J14").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,4,7,10);'

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("F15").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G15").Formula = '=A15&","&B15&","&C15&","&D15&")"'
$ws.Range("I15").Formula = '=F15&G15&";"'
$ws.Range("J15").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,3,7,11);'

$ws.Range("A16").Value = 2
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 12
$ws.Range("F16").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G16").Formula = '=A16&","&B16&","&C16&","&D16&")"'
$ws.Range("I16").Formula = '=F16&G16&";"'
$ws.Range("J16").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (2,2,7,12);'

$ws.Range("A18").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("F18").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G18").Formula = '=A18&","&B18&","&C18&","&D18&")"'
$ws.Range("I18").Formula = '=F18&G18&";"'
$ws.Range("J18").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,1,2,1);'

$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 2
$ws.Range("F19").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G19").Formula = '=A19&","&B19&","&C19&","&D19&")"'
$ws.Range("I19").Formula = '=F19&G19&";"'
$ws.Range("J19").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,2,3,2);'

$ws.Range("A20").Value = 1
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G20").Formula = '=A20&","&B20&","&C20&","&D20&")"'
$ws.Range("I20").Formula = '=F20&G20&";"'
$ws.Range("J20").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,3,4,3);'

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 4
$ws.Range("F21").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G21").Formula = '=A21&","&B21&","&C21&","&D21&")"'
$ws.Range("I21").Formula = '=F21&G21&";"'
$ws.Range("J21").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,4,5,4);'

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 5
$ws.Range("F22").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G22").Formula = '=A22&","&B22&","&C22&","&D22&")"'
$ws.Range("I22").Formula = '=F22&G22&";"'
$ws.Range("J22").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,4,6,5);'

$ws.Range("A23").Value = 1
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 6
$ws.Range("F23").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G23").Formula = '=A23&","&B23&","&C23&","&D23&")"'
$ws.Range("I23").Formula = '=F23&G23&";"'
$ws.Range("J23").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,2,4,6);'

$ws.Range("A24").Value = 1
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 7
$ws.Range("F24").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G24").Formula = '=A24&","&B24&","&C24&","&D24&")"'
$ws.Range("I24").Formula = '=F24&G24&";"'
$ws.Range("J24").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,1,7,7);'

$ws.Range("A25").Value = 1
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 8
$ws.Range("F25").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G25").Formula = '=A25&","&B25&","&C25&","&D25&")"'
$ws.Range("I25").Formula = '=F25&G25&";"'
$ws.Range("J25").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,5,7,8);'

$ws.Range("A26").Value = 1
$ws.Range("B26").Value = 5
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 9
$ws.Range("F26").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G26").Formula = '=A26&","&B26&","&C26&","&D26&")"'
$ws.Range("I26").Formula = '=F26&G26&";"'
$ws.Range("J26").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,5,4,9);'

$ws.Range("A27").Value = 1
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 10
$ws.Range("F27").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G27").Formula = '=A27&","&B27&","&C27&","&D27&")"'
$ws.Range("I27").Formula = '=F27&G27&";"'
$ws.Range("J27").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,4,7,10);'

$ws.Range("A28").Value = 1
$ws.Range("B28").Value = 3
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 11
$ws.Range("F28").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G28").Formula = '=A28&","&B28&","&C28&","&D28&")"'
$ws.Range("I28").Formula = '=F28&G28&";"'
$ws.Range("J28").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,3,7,11);'

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 12
$ws.Range("F29").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G29").Formula = '=A29&","&B29&","&C29&","&D29&")"'
$ws.Range("I29").Formula = '=F29&G29&";"'
$ws.Range("J29").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (1,2,7,12);'

$ws.Range("A31").Value = 3
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 1
$ws.Range("F31").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G31").Formula = '=A31&","&B31&","&C31&","&D31&")"'
$ws.Range("I31").Formula = '=F31&G31&";"'
$ws.Range("J31").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,1,2,1);'

$ws.Range("A32").Value = 3
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 2
$ws.Range("F32").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G32").Formula = '=A32&","&B32&","&C32&","&D32&")"'
$ws.Range("I32").Formula = '=F32&G32&";"'
$ws.Range("J32").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,2,3,2);'

$ws.Range("A33").Value = 3
$ws.Range("B33").Value = 3
$ws.Range("C33").Value = 4
$ws.Range("D33").Value = 3
$ws.Range("F33").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G33").Formula = '=A33&","&B33&","&C33&","&D33&")"'
$ws.Range("I33").Formula = '=F33&G33&";"'
$ws.Range("J33").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,3,4,3);'

$ws.Range("A34").Value = 3
$ws.Range("B34").Value = 4
$ws.Range("C34").Value = 5
$ws.Range("D34").Value = 4
$ws.Range("F34").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G34").Formula = '=A34&","&B34&","&C34&","&D34&")"'
$ws.Range("I34").Formula = '=F34&G34&";"'
$ws.Range("J34").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,4,5,4);'

$ws.Range("A35").Value = 3
$ws.Range("B35").Value = 4
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 5
$ws.Range("F35").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G35").Formula = '=A35&","&B35&","&C35&","&D35&")"'
$ws.Range("I35").Formula = '=F35&G35&";"'
$ws.Range("J35").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,4,6,5);'

$ws.Range("A36").Value = 3
$ws.Range("B36").Value = 2
$ws.Range("C36").Value = 4
$ws.Range("D36").Value = 6
$ws.Range("F36").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G36").Formula = '=A36&","&B36&","&C36&","&D36&")"'
$ws.Range("I36").Formula = '=F36&G36&";"'
$ws.Range("J36").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,2,4,6);'

$ws.Range("A37").Value = 3
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 7
$ws.Range("D37").Value = 7
$ws.Range("F37").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G37").Formula = '=A37&","&B37&","&C37&","&D37&")"'
$ws.Range("I37").Formula = '=F37&G37&";"'
$ws.Range("J37").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,1,7,7);'

$ws.Range("A38").Value = 3
$ws.Range("B38").Value = 5
$ws.Range("C38").Value = 7
$ws.Range("D38").Value = 8
$ws.Range("F38").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G38").Formula = '=A38&","&B38&","&C38&","&D38&")"'
$ws.Range("I38").Formula = '=F38&G38&";"'
$ws.Range("J38").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,5,7,8);'

$ws.Range("A39").Value = 3
$ws.Range("B39").Value = 5
$ws.Range("C39").Value = 4
$ws.Range("D39").Value = 9
$ws.Range("F39").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G39").Formula = '=A39&","&B39&","&C39&","&D39&")"'
$ws.Range("I39").Formula = '=F39&G39&";"'
$ws.Range("J39").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,5,4,9);'

$ws.Range("A40").Value = 3
$ws.Range("B40").Value = 4
$ws.Range("C40").Value = 7
$ws.Range("D40").Value = 10
$ws.Range("F40").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G40").Formula = '=A40&","&B40&","&C40&","&D40&")"'
$ws.Range("I40").Formula = '=F40&G40&";"'
$ws.Range("J40").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,4,7,10);'

$ws.Range("A41").Value = 3
$ws.Range("B41").Value = 3
$ws.Range("C41").Value = 7
$ws.Range("D41").Value = 11
$ws.Range("F41").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G41").Formula = '=A41&","&B41&","&C41&","&D41&")"'
$ws.Range("I41").Formula = '=F41&G41&";"'
$ws.Range("J41").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,3,7,11);'

$ws.Range("A42").Value = 3
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = 7
$ws.Range("D42").Value = 12
$ws.Range("F42").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values ('
$ws.Range("G42").Formula = '=A42&","&B42&","&C42&","&D42&")"'
$ws.Range("I42").Formula = '=F42&G42&";"'
$ws.Range("J42").Value = 'INSERT INTO sgr.flujo(id_workflow,id_estadoorigen,id_estadodestino,orden) values (3,2,7,12);'


# Apply a no-op style touch to the tail rows 1048560-1048576 of column J
# (mirrors the trailing artifact rows seen after formatting column J down
# to the sheet's full extent)
for ($r = 1048560; $r -le 1048576; $r++) {
    $ws.Cells.Item($r, 10).Font.Bold = $false
}


# Make "Flujo" the active sheet/tab with A1 selected (was "workflow" before)
$ws.Activate()
$ws.Range("A1").Select()
